$d = $word.ActiveDocument

$pairs = @(
    ,@("5+9=", "11+35=")
    ,@("28+9=", "25-10=")
    ,@("45+16=", "38-33=")
    ,@("5+88=", "21+14=")
    ,@("93+4=", "58+20=")
    ,@("5+52=", "59-52=")
    ,@("74-9=", "48-16=")
    ,@("22+15=", "74-29=")
    ,@("86-80=", "98-1=")
    ,@("12+21=", "22+56=")
    ,@("54-6=", "69+19=")
    ,@("28-11=", "75-29=")
    ,@("33-20=", "98-54=")
    ,@("82-75=", "5+35=")
    ,@("38-8=", "34+13=")
    ,@("34+24=", "37-14=")
    ,@("6+47=", "57-6=")
    ,@("15-8=", "25-20=")
    ,@("85-47=", "54-41=")
    ,@("90-61=", "33-11=")
    ,@("27-22=", "37+34=")
    ,@("67-32=", "48+23=")
    ,@("20+18=", "87-79=")
    ,@("11+16=", "44+31=")
    ,@("74-46=", "72-64=")
    ,@("6+14=", "21+61=")
    ,@("67-19=", "3+87=")
    ,@("62-16=", "69-31=")
    ,@("81-71=", "58+16=")
    ,@("48-36=", "0+83=")
    ,@("64-45=", "2+28=")
    ,@("6+4=", "75-12=")
    ,@("10+14=", "12+79=")
    ,@("37+62=", "15+63=")
    ,@("8+41=", "54-35=")
    ,@("95-38=", "23+76=")
    ,@("49+3=", "46+49=")
    ,@("89-36=", "33+6=")
    ,@("98-83=", "39+22=")
    ,@("30+23=", "13-5=")
    ,@("40+48=", "80+7=")
    ,@("93-84=", "26-23=")
    ,@("78+9=", "45+52=")
    ,@("98-24=", "85-76=")
    ,@("59-33=", "50+28=")
    ,@("95-93=", "71+28=")
    ,@("26+8=", "0+64=")
    ,@("23-14=", "76-71=")
    ,@("37+4=", "32+36=")
    ,@("86-0=", "97-48=")
    ,@("45-20=", "0+63=")
    ,@("38+54=", "26-26=")
    ,@("50+23=", "36+16=")
    ,@("25+35=", "67-40=")
    ,@("44+28=", "10+80=")
    ,@("32+10=", "28+45=")
    ,@("25+48=", "94-90=")
    ,@("17-12=", "79-50=")
    ,@("28+67=", "95-76=")
    ,@("12+49=", "53-21=")
    ,@("30+44=", "60+21=")
    ,@("22+55=", "48-32=")
    ,@("2+32=", "47+34=")
    ,@("13-9=", "50+9=")
    ,@("27+37=", "0+56=")
    ,@("69-22=", "8+31=")
    ,@("56+36=", "49+36=")
    ,@("77-63=", "44+24=")
    ,@("38+20=", "52-30=")
    ,@("83+10=", "65-62=")
    ,@("75+10=", "34+38=")
    ,@("63+6=", "50+48=")
    ,@("78-71=", "81-4=")
    ,@("1+58=", "13+47=")
    ,@("96-42=", "29+24=")
    ,@("89-72=", "54+17=")
    ,@("9+1=", "46+18=")
    ,@("40-18=", "72-9=")
    ,@("86+10=", "55-0=")
    ,@("18+39=", "26+37=")
    ,@("73-63=", "61+11=")
    ,@("90-43=", "76-67=")
    ,@("5-1=", "75+9=")
    ,@("56-10=", "66-20=")
    ,@("19+72=", "57+38=")
    ,@("66-35=", "36+26=")
    ,@("44+45=", "65+30=")
    ,@("45-33=", "30+31=")
    ,@("99-12=", "40+8=")
    ,@("2+87=", "77-51=")
    ,@("91-52=", "37+28=")
    ,@("28+63=", "54-4=")
    ,@("42+0=", "61-1=")
    ,@("16+54=", "84-41=")
    ,@("93-13=", "2+72=")
    ,@("79-68=", "41+0=")
    ,@("33+58=", "75+15=")
    ,@("70-55=", "70+14=")
    ,@("8+0=", "61-10=")
    ,@("97-92=", "75-66=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($pairs.Count) math problems"
